# Case_3_45 vm_pu results: update the 380 kV slack-bus voltage setpoint
# (B column, 1.05 -> 1.02 pu) and refresh the dependent bus-voltage
# results (columns C-F, I-N) for rows 2-25 to the recalculated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, followed by new values for columns B,C,D,E,F,I,J,K,L,M,N
$newValues = @(
    @(2, 1.02, 1.043199492457846, 1.050618818628143, 1.041229974745044, 1.061449140846215, 1.046747877627948, 1.04827122187715, 1.053372594034342, 1.044010079067248, 1.064173171589564, 1.019930378013375),
    @(3, 1.02, 1.044185484263137, 1.051401537821809, 1.042070524742142, 1.062396398459625, 1.047061685466904, 1.048903782298392, 1.053967872454583, 1.044661135540633, 1.064934718076539, 1.020145025167254),
    @(4, 1.02, 1.044823801540289, 1.051908272515548, 1.04261502520368, 1.063009967449704, 1.047263625999611, 1.049312792743105, 1.054352661503881, 1.045082375485458, 1.065427492437808, 1.020283708083621),
    @(5, 1.02, 1.045092224823668, 1.052121365672381, 1.042844077772312, 1.063268061571674, 1.047348254605782, 1.049484668549876, 1.054514331265698, 1.045259455053195, 1.065634654530064, 1.020341960213627),
    @(6, 1.02, 1.045137298636393, 1.052157148549043, 1.042882545153035, 1.063311405458109, 1.047362448439142, 1.049513522999868, 1.054541470704182, 1.045289186898994, 1.065669437915183, 1.020351738050892),
    @(7, 1.02, 1.044827387932589, 1.051911119634898, 1.042618085248583, 1.063013415529072, 1.04726475786166, 1.049315089640598, 1.054354822119583, 1.045084741670463, 1.065430260549968, 1.020284486648981),
    @(8, 1.02, 1.043532647476632, 1.050883287131607, 1.041513915928137, 1.061769139890949, 1.046854160859418, 1.048485059985659, 1.053573852668644, 1.044230113847177, 1.064430538635959, 1.020002961918797),
    @(9, 1.02, 1.041253586772192, 1.04907417705904, 1.039572935754195, 1.059581440482979, 1.046122124377978, 1.047020188336284, 1.052194686854385, 1.042723902039405, 1.062668963965138, 1.019505300108366),
    @(10, 1.02, 1.039735887035672, 1.047869553211842, 1.038282176921222, 1.058126314665449, 1.045628404246068, 1.046042134769491, 1.05127327249651, 1.04171964131404, 1.061494677581136, 1.019172483473592),
    @(11, 1.02, 1.03907911018075, 1.047348294556163, 1.03772404267801, 1.057497034501139, 1.045413272499515, 1.045618286131188, 1.050873831709008, 1.041284766665065, 1.060986232626342, 1.01902812706803),
    @(12, 1.02, 1.038835214151574, 1.047154729681907, 1.037516843886489, 1.057263412862123, 1.045333160914218, 1.045460798685907, 1.050725392852323, 1.041123232119205, 1.060797378699843, 1.018974470230751),
    @(13, 1.02, 1.038887527948783, 1.047196247553862, 1.037561283414666, 1.057313520003608, 1.045350354259176, 1.045494582583617, 1.050757236591537, 1.041157881939718, 1.060837888260752, 1.018985981448554),
    @(14, 1.02, 1.039058948429359, 1.04733229331998, 1.037706913166291, 1.057477720788733, 1.045406654570866, 1.045605269210357, 1.050861563105554, 1.041271414220713, 1.06097062179156, 1.019023692518717),
    @(15, 1.02, 1.039164574228049, 1.04741612276552, 1.037796656011928, 1.057578906379847, 1.045441316285937, 1.045673460130383, 1.050925833022965, 1.041341364865253, 1.061052404010699, 1.019046922734145),
    @(16, 1.02, 1.039779483559987, 1.047904154921677, 1.038319234821625, 1.058168094835426, 1.045642653471462, 1.046070256993885, 1.051299772385048, 1.041748502087801, 1.061528422076803, 1.019182058805557),
    @(17, 1.02, 1.040165306596543, 1.048210379392349, 1.038647242325916, 1.058537891788807, 1.045768586437901, 1.046319065271956, 1.051534211386496, 1.042003882838473, 1.061827023937191, 1.019266760822558),
    @(18, 1.02, 1.040390388933166, 1.048389028784815, 1.038838638055003, 1.058753665111706, 1.045841910949128, 1.04646415760141, 1.051670911038941, 1.042152839828699, 1.062001196072524, 1.019316142454347),
    @(19, 1.02, 1.040467142673189, 1.048449949320217, 1.03890391165864, 1.058827251271526, 1.045866890640553, 1.046513624663455, 1.051717514493264, 1.042203629903705, 1.062060584705803, 1.01933297629862),
    @(20, 1.02, 1.04012390747244, 1.048177520897255, 1.038612042511622, 1.05849820810196, 1.045755088471987, 1.046292373941184, 1.051509062913757, 1.041976483131401, 1.061794986495832, 1.019257675538807),
    @(21, 1.02, 1.039008467706469, 1.047292229739525, 1.037664025593073, 1.057429364400939, 1.045390081107408, 1.045572676165481, 1.050830843440088, 1.041237981863352, 1.060931534925596, 1.019012588551942),
    @(22, 1.02, 1.038307494629374, 1.046735923815626, 1.037068647832822, 1.056758040700437, 1.045159417539729, 1.045119877763774, 1.050404022439693, 1.040773641326941, 1.060388678972939, 1.018858281953704),
    @(23, 1.02, 1.038679060558741, 1.047030802224983, 1.037384204236171, 1.057113855336406, 1.045281807331149, 1.045359942591658, 1.050630325774909, 1.041019798251828, 1.060676454029121, 1.018940102668113),
    @(24, 1.02, 1.040142613836283, 1.048192368132328, 1.038627947562422, 1.058516139213978, 1.045761188023635, 1.046304434707652, 1.051520426564081, 1.041988863886709, 1.061809462828182, 1.019261780860471),
    @(25, 1.02, 1.04184248612804, 1.04954162450547, 1.040074161325043, 1.06014642912941, 1.04631237928422, 1.047399155424853, 1.052551584894354, 1.043113318442328, 1.063124360175963, 1.019634142735887)
)

$cols = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

foreach ($entry in $newValues) {
    $r = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $entry[$i + 1]
    }
}
